$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $rng = $ws.Range($Cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = $origStyle
}

Set-TextValue "D2" "326.39"
Set-TextValue "E2" "-0.10%"
Set-TextValue "D3" "44.24"
Set-TextValue "E3" "-1.66%"
Set-TextValue "D4" "5.509"
Set-TextValue "E4" "-0.96%"
Set-TextValue "D5" "0.08024"
Set-TextValue "E5" "-0.96%"
Set-TextValue "D6" "1.996"
Set-TextValue "E6" "4.61%"
Set-TextValue "D7" "4.293"
Set-TextValue "E7" "-1.00%"
Set-TextValue "D8" "2.572"
Set-TextValue "E8" "-6.17%"
Set-TextValue "D9" "0.9488"
Set-TextValue "E9" "-0.23%"
Set-TextValue "D10" "0.1145"
Set-TextValue "E10" "-1.45%"
Set-TextValue "D11" "0.1836"
Set-TextValue "E11" "-3.35%"
Set-TextValue "D12" "12.18"
Set-TextValue "E12" "42.76%"
Set-TextValue "D13" "0.09743"
Set-TextValue "E13" "-5.10%"
Set-TextValue "D14" "0.04616"
Set-TextValue "E14" "11.28%"
Set-TextValue "E15" "0.00%"
Set-TextValue "D16" "0.001273"
Set-TextValue "E16" "0.08%"
Set-TextValue "E17" "-4.51%"
Set-TextValue "D18" "0.005871"
Set-TextValue "E18" "-3.76%"
Set-TextValue "B19" "HotbitToken"
Set-TextValue "C19" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D19" "0.004301"
Set-TextValue "E19" "-6.21%"
Set-TextValue "B20" "LEO"
Set-TextValue "C20" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D20" "3.367"
Set-TextValue "E20" "-6.69%"
Set-TextValue "B21" "BitpandaEcosystemToken"
Set-TextValue "C21" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D21" "0.3479"
Set-TextValue "E21" "-0.20%"
Set-TextValue "B22" "ProBitToken"
Set-TextValue "C22" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D22" "0.1409"
Set-TextValue "E22" "2.56%"
Set-TextValue "B23" "ZBToken"
Set-TextValue "C23" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D23" "0.2545"
Set-TextValue "E23" "-4.51%"
Set-TextValue "B24" "BitKan"
Set-TextValue "C24" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D24" "0.001244"
Set-TextValue "E24" "0.28%"
Set-TextValue "D25" "0.0001189"
Set-TextValue "E25" "-3.68%"
Set-TextValue "E26" "-6.50%"
Set-TextValue "D38" "0.02556"
Set-TextValue "E38" "-4.53%"
Set-TextValue "D39" "0.05531"
Set-TextValue "E39" "-0.96%"
Set-TextValue "D40" "0.007515"
Set-TextValue "E40" "-1.69%"
Set-TextValue "D41" "0.1392"
Set-TextValue "E41" "-0.23%"
Set-TextValue "D42" "0.007612"
Set-TextValue "E42" "-32.84%"
Set-TextValue "D43" "0.002013"
Set-TextValue "E43" "-2.30%"
Set-TextValue "D44" "0.008511"
Set-TextValue "E44" "-2.09%"
Set-TextValue "D45" "0.00007113"
Set-TextValue "E45" "0.21%"
Set-TextValue "E46" "-0.39%"
Set-TextValue "E47" "1.07%"
Set-TextValue "D48" "0.003357"
Set-TextValue "E48" "-3.08%"
Set-TextValue "E49" "-0.39%"
Set-TextValue "E50" "-0.39%"
